# The source notebook/script was restructured to push the scraped rows
# through a pandas DataFrame (and stop saving each trademark image to
# disk) before writing the workbook back out. That rewrite changed two
# observable things in the workbook itself:
#
#   1. The sheet got renamed from the default "Sheet1" to "Sheet" (the
#      name pandas' `ExcelWriter` uses unless told otherwise).
#   2. The bespoke header formatting (bold font, thin box border,
#      centered/top-aligned text) that had been applied to row 1 is gone
#      - a plain `DataFrame.to_excel()` writes header cells with no
#      special styling at all, so the custom style is removed.
#
# Apply both changes here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet: "Sheet1" -> "Sheet"
$ws.Name = "Sheet"

# 2) Strip the header row's custom formatting (bold / border / alignment)
#    so it falls back to the workbook's plain default "Normal" style,
#    matching a freshly-written DataFrame export.
$header = $ws.Range("A1:F1")
$header.Style = "Normal"
